$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the trailing placeholder row (34) that only held leftover blue/bold
# formatting (K34:O34) with no data.
$ws.Rows(34).Delete()

# Append the new data row (26) for 27/03/2020 ("Se agrega excel de decesos").
$ws.Range("A26").Value = 43917
$ws.Range("A26").NumberFormat = "DD/MM/YY"

$ws.Range("B26").Value = 25
$ws.Range("C26").Value = 3
$ws.Range("D26").Value = 5
$ws.Range("E26").Value = 25
$ws.Range("F26").Value = 1
$ws.Range("G26").Value = 15
$ws.Range("H26").Value = 71
$ws.Range("I26").Value = 1084
$ws.Range("J26").Value = 20
$ws.Range("K26").Value = 37
$ws.Range("L26").Value = 172
$ws.Range("M26").Value = 151
$ws.Range("N26").Value = 177
$ws.Range("O26").Value = 27
$ws.Range("P26").Value = 93
$ws.Range("Q26").Value = 2
$ws.Range("R26").Value = 26
$ws.Range("S26").Value = 1909

# Move/restore the selection to just past the new last row, like the source
# workbook shows after the edit.
[void]$ws.Range("A27").Select()
